$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("BQ15").Select()
$excel.ActiveWindow.ScrollColumn = 68
$excel.ActiveWindow.ScrollRow = 11
